$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '34.431.62'
$ws.Range('E2').Value = '  +0.31%  '

# Row 3
$ws.Range('D3').Value = '1.799.05'
$ws.Range('E3').Value = '  +0.28%  '

# Row 4
$ws.Range('E4').Value = '  +0.33%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '224.28'
$ws.Range('E5').Value = '  -0.37%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.601'
$ws.Range('E6').Value = '  +1.70%  '

# Row 7
$ws.Range('E7').Value = '  +0.22%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.55'
$ws.Range('E8').Value = '  +15.25%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.289'
$ws.Range('E9').Value = '  -0.12%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0665'
$ws.Range('E10').Value = '  -1.26%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0995'
$ws.Range('E11').Value = '  +3.50%  '

# Row 12
$ws.Range('D12').Value = '2.058.67'
$ws.Range('E12').Value = '  +0.26%  '

# Row 13
$ws.Range('D13').Value = '1.798.62'
$ws.Range('E13').Value = '  -0.06%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.82'
$ws.Range('E14').Value = '  -2.70%  '

# Row 15
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').Value = '34.445.92'
$ws.Range('E15').Value = '  +0.37%  '

# Row 16
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.624'
$ws.Range('E16').Value = '  -0.45%  '

# Row 17
$ws.Range('E17').Value = '  +0.19%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '67.24'
$ws.Range('E18').Value = '  -2.15%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '239.73'
$ws.Range('E19').Value = '  -0.10%  '

# Row 20
$ws.Range('D20').Value = '0.0₃0765'
$ws.Range('E20').Value = '  -0.41%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.08'
$ws.Range('E21').Value = '  -1.08%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  +0.25%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.07'
$ws.Range('E23').Value = '  -0.15%  '

# Row 24
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.15'
$ws.Range('E24').Value = '  -1.76%  '

# Row 25
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.51'
$ws.Range('E25').Value = '  +1.10%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.62'
$ws.Range('E26').Value = '  -3.17%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.30'
$ws.Range('E27').Value = '  +0.92%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.120'
$ws.Range('E28').Value = '  +0.78%  '

# Row 29
$ws.Range('E29').Value = '  +0.33%  '

# Row 30
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.23'
$ws.Range('E30').Value = '  -0.14%  '

# Row 31
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.77'
$ws.Range('E31').Value = '  +0.32%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0510'
$ws.Range('E32').Value = '  +0.23%  '

# Row 33
$ws.Range('E33').Value = '  -1.00%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.76'
$ws.Range('E34').Value = '  +0.65%  '

# Row 35
$ws.Range('B35').Value = 'Maker'
$ws.Range('C35').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D35').Value = '1.314.76'
$ws.Range('E35').Value = '  -2.94%  '

# Row 36
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.642'
$ws.Range('E36').Value = '  +0.22%  '

# Row 37
$ws.Range('E37').Value = '  +0.72%  '

# Row 38
$ws.Range('B38').Value = 'Aave'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '85.68'
$ws.Range('E38').Value = '  +6.51%  '

# Row 39
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0187'
$ws.Range('E39').Value = '  +1.76%  '

# Row 40
$ws.Range('B40').Value = 'InjectiveProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '14.67'
$ws.Range('E40').Value = '  +12.30%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.24'
$ws.Range('E41').Value = '  +5.91%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.31'
$ws.Range('E42').Value = '  -0.58%  '

# Row 43
$ws.Range('E43').Value = '  +0.86%  '

# Row 44
$ws.Range('E44').Value = '  +0.30%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.937'
$ws.Range('E45').Value = '  +1.10%  '

# Row 46
$ws.Range('E46').Value = '  +4.94%  '

# Row 47
$ws.Range('D47').Value = '1.959.58'
$ws.Range('E47').Value = '  +0.27%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.83'
$ws.Range('E48').Value = '  +1.39%  '

# Row 49
$ws.Range('E49').Value = '  +0.28%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '100.51'
$ws.Range('E50').Value = '  -0.91%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0609'
$ws.Range('E51').Value = '  +1.42%  '
